$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsOverview.Range("D7").Value = "2016-34-18 20:34:03"
$wsZhCn.Range("E7").Value = "2016-03-18 20:33:59"
$wsDeDe.Range("E7").Value = "2016-03-18 20:34:03"
